# Dev #621 - remove hardcoding
# Rename the "Tab" worksheet to "KS4" (defined names / autofilter refs
# that point at the sheet follow the rename automatically), and add the
# underlying formula to the existing "contains text" conditional format
# rule on column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab")
$ws.Name = "KS4"

$rng = $ws.Range("L7:L114,L123:L1243")
$fc = $rng.FormatConditions.Item(1)
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Yes",L7)))'
